# Revisi / Tambah Fitur
#
# Adds the six new worksheets (Gejala, Pasien, Penyakit, Rule, gejala_pasien,
# hasil_diagnosa) that make up the rest of the "Raw Database" schema, next to
# the pre-existing "Users" sheet.
#
# Cell values are written in the same chronological order the original
# author used (verified against the shared-string append order in the
# target file) so that the shared string table ends up byte-identical:
# both tables' "kode_penyakit" column was filled in last, after every other
# sheet/column had already been populated.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Create all six new sheets first (tab order: Users, Gejala, Pasien,
# Penyakit, Rule, gejala_pasien, hasil_diagnosa)
# ---------------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$gejala = $wb.Worksheets.Add($null, $last)
$gejala.Name = "Gejala"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$pasien = $wb.Worksheets.Add($null, $last)
$pasien.Name = "Pasien"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$penyakit = $wb.Worksheets.Add($null, $last)
$penyakit.Name = "Penyakit"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$rule = $wb.Worksheets.Add($null, $last)
$rule.Name = "Rule"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$gejalaPasien = $wb.Worksheets.Add($null, $last)
$gejalaPasien.Name = "gejala_pasien"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$hasil = $wb.Worksheets.Add($null, $last)
$hasil.Name = "hasil_diagnosa"

# ---------------------------------------------------------------------------
# Gejala
# ---------------------------------------------------------------------------
$gejala.Range("A1").Value = "id_gejala"
$gejala.Range("A2").Value = "kode_gejala"
$gejala.Range("A3").Value = "nama_gejala"
$gejala.Range("A4").Value = "created_at"
$gejala.Range("A5").Value = "updated_at"

# ---------------------------------------------------------------------------
# Pasien
# ---------------------------------------------------------------------------
$pasien.Range("A1").Value = "id_pasien"
$pasien.Range("A2").Value = "user_id"
$pasien.Range("A3").Value = "nama_pasien"
$pasien.Range("A4").Value = "nik"
$pasien.Range("A5").Value = "jenis_kelamin"
$pasien.Range("A6").Value = "tanggal_lahir"
$pasien.Range("A7").Value = "usia"
$pasien.Range("A8").Value = "alamat"
$pasien.Range("A9").Value = "no_hp"
$pasien.Range("A10").Value = "created_at"
$pasien.Range("A11").Value = "updated_at"

# ---------------------------------------------------------------------------
# Penyakit -- "kode_penyakit" (row 2) is filled in at the very end, below
# ---------------------------------------------------------------------------
$penyakit.Range("A1").Value = "id_penyakit"
$penyakit.Range("A3").Value = "nama_penyakit"
$penyakit.Range("A4").Value = "deskripsi"
$penyakit.Range("A5").Value = "solusi"

# ---------------------------------------------------------------------------
# Rule
# ---------------------------------------------------------------------------
$rule.Range("A1").Value = "id_aturan"
$rule.Range("A2").Value = "id_penyakit"
$rule.Range("A3").Value = "id_gejala"
$rule.Range("A5").Value = "nilai_md"
$rule.Range("A4").Value = "nilai_mb"

# ---------------------------------------------------------------------------
# gejala_pasien
# ---------------------------------------------------------------------------
$gejalaPasien.Range("A1").Value = "id_gejala_pasien"
$gejalaPasien.Range("A2").Value = "id_pasien"
$gejalaPasien.Range("A3").Value = "id_gejala"
$gejalaPasien.Range("A4").Value = "cf_user"

# ---------------------------------------------------------------------------
# hasil_diagnosa -- "kode_penyakit" (row 4) is filled in at the very end, below
# ---------------------------------------------------------------------------
$hasil.Range("A1").Value = "id_hasil"
$hasil.Range("A2").Value = "id_pasien"
$hasil.Range("A3").Value = "nama_pasien"
$hasil.Range("A5").Value = "tanggal_diagnosa"
$hasil.Range("A6").Value = "nilai_cf"
$hasil.Range("A7").Value = "diagnosa"
$hasil.Range("A8").Value = "keterangan"

# ---------------------------------------------------------------------------
# Go back and add the "kode_penyakit" column that was missing from both
# Penyakit and hasil_diagnosa
# ---------------------------------------------------------------------------
$penyakit.Range("A2").Value = "kode_penyakit"
$hasil.Range("A4").Value = "kode_penyakit"

# ---------------------------------------------------------------------------
# Column widths (best-fit-ish) and print setup
# ---------------------------------------------------------------------------
$pasien.Columns.Item(1).ColumnWidth = 12.6
$pasien.PageSetup.Orientation = 1
$gejalaPasien.Columns.Item(1).ColumnWidth = 15.15
$hasil.Columns.Item(1).ColumnWidth = 15.65

# ---------------------------------------------------------------------------
# Selections on each sheet
# ---------------------------------------------------------------------------
[void]$gejala.Range("A6").Select()
[void]$pasien.Range("J16").Select()
[void]$penyakit.Range("G6").Select()
[void]$rule.Range("A6").Select()
[void]$gejalaPasien.Range("H7").Select()
[void]$hasil.Range("K26").Select()

# Leave the last sheet (hasil_diagnosa) as the active / selected tab
$hasil.Activate()
